$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing weekly data rows ---
$ws.Range("B40").Value = 764
$ws.Range("B41").Value = 823
$ws.Range("B42").Value = 997

# Row 45 will become the new "most recent" row, so grab the red-highlighted
# date style (currently on A42) for it before A42's style is changed below.
$ws.Range("A42").Copy()
$ws.Range("A45").PasteSpecial(-4122)

# Row 42 is no longer the latest week once rows 43-45 are appended, so its
# special red-highlighted date style is replaced with the normal highlighted
# date style (still highlighted fill, default font color).
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").NumberFormat = "m/d/yyyy"
$ws.Range("A42").Font.ThemeColor = 1

# --- New weekly data rows 43-45 ---
$ws.Range("A43").Value = 44191
$ws.Range("B43").Value = 975
$ws.Range("C42").Copy()
$ws.Range("C43").PasteSpecial(-4122)
$ws.Range("C43").Formula = "=(B43-B42)/B43"
$ws.Range("D43").Formula = "=(B43/202558)*100000"

$ws.Range("A44").Value = 44198
$ws.Range("B44").Value = 1104
$ws.Range("C42").Copy()
$ws.Range("C44").PasteSpecial(-4122)
$ws.Range("C44").Formula = "=(B44-B43)/B44"
$ws.Range("D44").Formula = "=(B44/202558)*100000"

$ws.Range("A45").Value = 44205
$ws.Range("B45").Value = 745
$ws.Range("C42").Copy()
$ws.Range("C45").PasteSpecial(-4122)
$ws.Range("C45").Formula = "=(B45-B44)/B45"
$ws.Range("D45").Formula = "=(B45/202558)*100000"

# --- Update selection to reflect the new last-entered cell ---
[void]$ws.Range("B44").Select()
